$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IHW")

# --- Insert a new configuration row at row 64 (everything below shifts down by one) ---
$ws.Rows.Item(64).Insert()

$ws.Cells.Item(64, 1).Value = "includemodelinfo"
$ws.Cells.Item(64, 2).Value = "yes|no"
$ws.Cells.Item(64, 3).Clear()
$ws.Cells.Item(64, 4).Value = $false
$ws.Cells.Item(64, 5).Value = $true
$ws.Cells.Item(64, 6).Value = "no"

# --- View changes: zoom in, scroll the frozen pane, move the active selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 175
$excel.ActiveWindow.ScrollRow = 59
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("G64").Select()
